# Generate Report for Handoff
# Regenerate the localization handoff files (new GUID-based names / hashes)
# and refresh the report: the Overview sheet and the per-locale (zh-cn,
# de-de) sheets all reference the old source file "d88e5db0-..." which is
# replaced by the freshly generated "c4e7e47a-...". The per-locale sheets
# also get a brand new "Latest Handoff File" / "Latest Handoff Datetime",
# while the "Latest Target File" / "Latest Handback File" (and their
# hyperlinks) are cleared out and "Latest Handback DateTime" is reset to
# the zero date, since the new handoff hasn't been handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "d88e5db0-5aa9-468e-87d8-9526b2ea6769"
$newGuid = "c4e7e47a-b530-4695-8b5c-9ddc9d28607b"

$oldHash = "2e5cee938237a608871ae136aecc6501d2ac445d"
$newHash = "8c7977dc8acac4bf0ae054086d4096bff951bcfe"

$newFileName = $newGuid + ".md"
$newPathName = "e2e\" + $newGuid + ".md"

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = $newFileName
$ov.Range("B2").Value = $newPathName
$ov.Hyperlinks.Item(1).TextToDisplay = $newPathName
$ov.Range("G2").Value = "2016-08-26 00:59:14"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = $newFileName
$zh.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-26 00:59:09"
$zh.Hyperlinks.Item(1).TextToDisplay = $newFileName
$zh.Hyperlinks.Item(2).Delete()
$zh.Range("I2").Style = "Normal"
$zh.Range("I2").Value = ""
$zh.Range("J2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = $newFileName
$de.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$de.Range("H2").Value = "2016-08-26 00:59:14"
$de.Hyperlinks.Item(1).TextToDisplay = $newFileName
$de.Hyperlinks.Item(2).Delete()
$de.Range("I2").Style = "Normal"
$de.Range("I2").Value = ""
$de.Range("J2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426
